$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("day")

# --- Convert column D (bsecode) in rows 632-644 from text to numeric ---
foreach ($r in 632..644) {
    $ws.Cells.Item($r, 4).Value = $ws.Cells.Item($r, 4).Value() + 0
}

# --- Append new data rows 645-657 (next-day "02/10/2024" snapshot) ---
# row 645
$ws.Cells.Item(645, 1).Value = 1
$ws.Cells.Item(645, 2).Value = "OFSS"
$ws.Cells.Item(645, 3).Value = "Oracle Financial Services Software Limited"
$c = $ws.Cells.Item(645, 4)
$c.NumberFormat = "@"
$c.Value = "532466"
$c.Style = "Normal"
$ws.Cells.Item(645, 5).Value = -0.14
$ws.Cells.Item(645, 6).Value = 11442.65
$ws.Cells.Item(645, 7).Value = 54125
$ws.Cells.Item(645, 8).Value = "day"
$ws.Cells.Item(645, 9).Value = "02/10/2024 11:35:36"

# row 646
$ws.Cells.Item(646, 1).Value = 2
$ws.Cells.Item(646, 2).Value = "MPHASIS"
$ws.Cells.Item(646, 3).Value = "Mphasis Limited"
$c = $ws.Cells.Item(646, 4)
$c.NumberFormat = "@"
$c.Value = "526299"
$c.Style = "Normal"
$ws.Cells.Item(646, 5).Value = 0.17
$ws.Cells.Item(646, 6).Value = 3015.6
$ws.Cells.Item(646, 7).Value = 899959
$ws.Cells.Item(646, 8).Value = "day"
$ws.Cells.Item(646, 9).Value = "02/10/2024 11:35:36"

# row 647
$ws.Cells.Item(647, 1).Value = 3
$ws.Cells.Item(647, 2).Value = "ICICIGI"
$ws.Cells.Item(647, 3).Value = "ICICI Lombard General Insurance Company Ltd"
$c = $ws.Cells.Item(647, 4)
$c.NumberFormat = "@"
$c.Value = "540716"
$c.Style = "Normal"
$ws.Cells.Item(647, 5).Value = -0.9399999999999999
$ws.Cells.Item(647, 6).Value = 2153.9
$ws.Cells.Item(647, 7).Value = 381613
$ws.Cells.Item(647, 8).Value = "day"
$ws.Cells.Item(647, 9).Value = "02/10/2024 11:35:36"

# row 648
$ws.Cells.Item(648, 1).Value = 4
$ws.Cells.Item(648, 2).Value = "SBILIFE"
$ws.Cells.Item(648, 3).Value = "SBI Life Insurance Company Ltd"
$c = $ws.Cells.Item(648, 4)
$c.NumberFormat = "@"
$c.Value = "540719"
$c.Style = "Normal"
$ws.Cells.Item(648, 5).Value = -0.53
$ws.Cells.Item(648, 6).Value = 1834.2
$ws.Cells.Item(648, 7).Value = 885515
$ws.Cells.Item(648, 8).Value = "day"
$ws.Cells.Item(648, 9).Value = "02/10/2024 11:35:36"

# row 649
$ws.Cells.Item(649, 1).Value = 5
$ws.Cells.Item(649, 2).Value = "CIPLA"
$ws.Cells.Item(649, 3).Value = "Cipla Limited"
$c = $ws.Cells.Item(649, 4)
$c.NumberFormat = "@"
$c.Value = "500087"
$c.Style = "Normal"
$ws.Cells.Item(649, 5).Value = 0.65
$ws.Cells.Item(649, 6).Value = 1664.85
$ws.Cells.Item(649, 7).Value = 945740
$ws.Cells.Item(649, 8).Value = "day"
$ws.Cells.Item(649, 9).Value = "02/10/2024 11:35:36"

# row 650
$ws.Cells.Item(650, 1).Value = 6
$ws.Cells.Item(650, 2).Value = "TECHM"
$ws.Cells.Item(650, 3).Value = "Tech Mahindra Limited"
$c = $ws.Cells.Item(650, 4)
$c.NumberFormat = "@"
$c.Value = "532755"
$c.Style = "Normal"
$ws.Cells.Item(650, 5).Value = 3.06
$ws.Cells.Item(650, 6).Value = 1625.4
$ws.Cells.Item(650, 7).Value = 5176596
$ws.Cells.Item(650, 8).Value = "day"
$ws.Cells.Item(650, 9).Value = "02/10/2024 11:35:36"

# row 651
$ws.Cells.Item(651, 1).Value = 7
$ws.Cells.Item(651, 2).Value = "TATACONSUM"
$ws.Cells.Item(651, 3).Value = "TATA Consumer Products Ltd"
$c = $ws.Cells.Item(651, 4)
$c.NumberFormat = "@"
$c.Value = "500800"
$c.Style = "Normal"
$ws.Cells.Item(651, 5).Value = -0.06
$ws.Cells.Item(651, 6).Value = 1196.25
$ws.Cells.Item(651, 7).Value = 892381
$ws.Cells.Item(651, 8).Value = "day"
$ws.Cells.Item(651, 9).Value = "02/10/2024 11:35:36"

# row 652
$ws.Cells.Item(652, 1).Value = 8
$ws.Cells.Item(652, 2).Value = "JINDALSTEL"
$ws.Cells.Item(652, 3).Value = "Jindal Steel & Power Limited"
$c = $ws.Cells.Item(652, 4)
$c.NumberFormat = "@"
$c.Value = "532286"
$c.Style = "Normal"
$ws.Cells.Item(652, 5).Value = -0.41
$ws.Cells.Item(652, 6).Value = 1035.35
$ws.Cells.Item(652, 7).Value = 2382544
$ws.Cells.Item(652, 8).Value = "day"
$ws.Cells.Item(652, 9).Value = "02/10/2024 11:35:36"

# row 653
$ws.Cells.Item(653, 1).Value = 9
$ws.Cells.Item(653, 2).Value = "HDFCLIFE"
$ws.Cells.Item(653, 3).Value = "HDFC Life Insurance Company Ltd"
$c = $ws.Cells.Item(653, 4)
$c.NumberFormat = "@"
$c.Value = "540777"
$c.Style = "Normal"
$ws.Cells.Item(653, 5).Value = -1.09
$ws.Cells.Item(653, 6).Value = 710.2
$ws.Cells.Item(653, 7).Value = 1987048
$ws.Cells.Item(653, 8).Value = "day"
$ws.Cells.Item(653, 9).Value = "02/10/2024 11:35:36"

# row 654
$ws.Cells.Item(654, 1).Value = 10
$ws.Cells.Item(654, 2).Value = "BERGEPAINT"
$ws.Cells.Item(654, 3).Value = "Berger Paints (i) Limited"
$c = $ws.Cells.Item(654, 4)
$c.NumberFormat = "@"
$c.Value = "509480"
$c.Style = "Normal"
$ws.Cells.Item(654, 5).Value = -0.43
$ws.Cells.Item(654, 6).Value = 619.65
$ws.Cells.Item(654, 7).Value = 1098938
$ws.Cells.Item(654, 8).Value = "day"
$ws.Cells.Item(654, 9).Value = "02/10/2024 11:35:36"

# row 655
$ws.Cells.Item(655, 1).Value = 11
$ws.Cells.Item(655, 2).Value = "DABUR"
$ws.Cells.Item(655, 3).Value = "Dabur India Limited"
$c = $ws.Cells.Item(655, 4)
$c.NumberFormat = "@"
$c.Value = "500096"
$c.Style = "Normal"
$ws.Cells.Item(655, 5).Value = -0.99
$ws.Cells.Item(655, 6).Value = 619
$ws.Cells.Item(655, 7).Value = 1095384
$ws.Cells.Item(655, 8).Value = "day"
$ws.Cells.Item(655, 9).Value = "02/10/2024 11:35:36"

# row 656
$ws.Cells.Item(656, 1).Value = 12
$ws.Cells.Item(656, 2).Value = "BIOCON"
$ws.Cells.Item(656, 3).Value = "Biocon Limited"
$c = $ws.Cells.Item(656, 4)
$c.NumberFormat = "@"
$c.Value = "532523"
$c.Style = "Normal"
$ws.Cells.Item(656, 5).Value = 2.08
$ws.Cells.Item(656, 6).Value = 370.65
$ws.Cells.Item(656, 7).Value = 2195117
$ws.Cells.Item(656, 8).Value = "day"
$ws.Cells.Item(656, 9).Value = "02/10/2024 11:35:36"

# row 657
$ws.Cells.Item(657, 1).Value = 13
$ws.Cells.Item(657, 2).Value = "GMRINFRA"
$ws.Cells.Item(657, 3).Value = "Gmr Infrastructure Limited"
$c = $ws.Cells.Item(657, 4)
$c.NumberFormat = "@"
$c.Value = "532754"
$c.Style = "Normal"
$ws.Cells.Item(657, 5).Value = -0.14
$ws.Cells.Item(657, 6).Value = 93.93000000000001
$ws.Cells.Item(657, 7).Value = 6427524
$ws.Cells.Item(657, 8).Value = "day"
$ws.Cells.Item(657, 9).Value = "02/10/2024 11:35:36"
Write-Output "done"
